# skoleeksamen.xlsx: "smaaplukk og v25 eksamen lagt til paa ny"
# - Column D header ("Ikke lenger pensum per ...") moved from the 2023 note
#   to a 2025 note, since the 2023 one is no longer needed (2023 row content
#   unchanged, only the footnote text shifted forward to 2025).
# - New column E ("Data") added, with a link to the extra material for the
#   2025 - Vaar exam in the last data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D1: the "no longer on syllabus" note now refers to 2025 (was 2023)
$ws.Range("D1").Value = "Ikke lenger pensum per 2025"

# New "Data" column header
$ws.Range("E1").Value = "Data"

# New extra-material link for the 2025 - Vaar exam row
$ws.Range("E18").Value = "[Materiale](tidligere-eksamensoppgaver/skole-25-v-ekstra.zip)"

# Give the new column a sensible width (best effort autosize, matches the
# other "link" columns' treatment) and move selection to reflect where the
# editor ended up after adding the new column/row content.
$ws.Columns.Item(5).AutoFit() | Out-Null

$ws.Range("E19").Select() | Out-Null
